$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# for every data row (row 2 through row 454) on each automatic update run.
$lastRow = 454
$ws.Range("C2:C$lastRow").Value = 45182
